$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Matthew Darby's logged time increased from 60h 35m to 66h 20m (+5:45)
$ws.Range("B4").Value = "66h 20m"

# Move the active selection to B4, matching the edited workbook state
$ws.Range("B4").Select()
